$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "23.446.17"
$ws.Range("E2").Value = "  +1.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.638.89"
$ws.Range("E3").Value = "  +2.36%  "

# Row 4 - TetherUSD
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - USDC
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

# Row 6 - BNB
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "306.31"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.03%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.55%  "

# Row 8 - OKB
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "52.20"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "

# Row 9 - Cardano
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.3644"
$cell.Style = "Normal"

# Row 10 - Polygon
$ws.Range("E10").Value = "  -0.40%  "

# Row 11 - Dogecoin
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.08148"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "

# Row 12 - BinanceUSD
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "

# Row 13 - Solana
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "22.96"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.24%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.85%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +2.60%  "

# Row 16 - Chainlink
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "7.373"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.638.25"
$ws.Range("E17").Value = "  +2.22%  "

# Row 18 - Litecoin
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "94.77"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19 - TRON
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06922"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "

# Row 20 - Avalanche
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "18.18"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "23.459.05"
$ws.Range("E23").Value = "  +1.19%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -1.36%  "

# Row 25 - LidoDAOToken
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.072"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +2.84%  "

# Row 26 - Toncoin
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.417"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.70%  "

# Row 27 - EthereumClassic
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "21.23"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

# Row 28 - Monero
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "150.77"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.89%  "

# Row 29 - HuobiToken
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.348"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "

# Row 30 - BitcoinCash
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "137.62"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.75%  "

# Row 31 - WEMIXTOKEN
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.295"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -3.66%  "

# Row 32 - WrappedliquidstakedEther2.0
$ws.Range("D32").Value = "1.818.53"
$ws.Range("E32").Value = "  +2.20%  "

# Row 33 - Filecoin
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.798"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "

# Row 34 - ImmutableX
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.9661"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "

# Row 35 - VeChain
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.02844"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +4.97%  "

# Row 36 - FraxShare
$ws.Range("E36").Value = "  +0.73%  "

# Row 37 - Hedera
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.07320"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.28%  "

# Row 38 - Algorand
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2530"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.17%  "

# Row 41 - TrustWalletToken
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.378"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "

# Row 42 - TheSandbox
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.7100"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43 - Aptos
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "12.50"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "

# Row 44 - EnergySwap
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.14"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.12%  "

# Row 45 - Decentraland
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.6555"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

# Row 46 - NEARProtocol
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.339"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.21%  "

# Row 47 - Frax
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "

# Row 48 - PancakeSwap
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "4.018"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49 - Cronos
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.07966"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

# Row 50 - Quant
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "129.00"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.29%  "

# Row 51 - Flow
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.204"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "

# Row 39 & 40 - swap Stellar and InternetComputer(DFINITY), with new data
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.132"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.08833"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
